$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) text updates ---
$ws.Range("N1").Value = "P_charge_lambda"
$ws.Range("O1").Value = "v"
$ws.Range("P1").Value = "v_variance"
$ws.Range("Q1").Value = "SOC_warn"

# --- Row 2 data updates ---
$ws.Range("G2").Value = 28
$ws.Range("H2").Value = 15
$ws.Range("J2").Value = 72
$ws.Range("K2").Value = 30
$ws.Range("N2").Value = 0.9
$ws.Range("O2").Value = 1.1
$ws.Range("P2").Value = 0.1
$ws.Range("Q2").Value = 0.4

# --- Row 3 data updates ---
$ws.Range("G3").Value = 28
$ws.Range("H3").Value = 15
$ws.Range("J3").Value = 72
$ws.Range("K3").Value = 30
$ws.Range("N3").Value = 0.9
$ws.Range("O3").Value = 1.1
$ws.Range("P3").Value = 0.1
$ws.Range("Q3").Value = 0.4

# --- Column width: extend the 14-width group to include column 20 ---
$ws.Range($ws.Cells.Item(1, 20), $ws.Cells.Item(1, 20)).ColumnWidth = 14

# --- Sheet view: scroll + selection ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("K3").Select()

# --- Window size/position ---
$excel.ActiveWindow.Left = -120
$excel.ActiveWindow.Top = -120
$excel.ActiveWindow.Width = 29040
$excel.ActiveWindow.Height = 15840
